# Apply the Wellness.xlsx update: append 16 new rows (126-141) of player wellness data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formatting (styles) from the last existing data row (125) down over the new rows (126:141)
$ws.Range("A125:I125").Copy()
$ws.Range("A126:I141").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Data for the new rows: Row, Date, Player, Volume, Intensite, Charge, Fatigue, Localisation, Plaisir
$newRows = @(
    @(126, 45874, "Amir Etien", 70, 6, 6, 7, "Ischios", 10),
    @(127, 45874, "Omar Benyounes", 70, 6, 7, 0, "", 6),
    @(128, 45874, "Romain Thunet", 70, 8, 7, 0, "", 7),
    @(129, 45874, "Jeremie Laurent", 70, 8, 7, 3, "Courbatures ", 8),
    @(130, 45874, "Emmanuel Valey", 70, 8, 7, 5, "Adducteur", 3),
    @(131, 45874, "Amir Kherrab", 70, 7, 7, 0, "", 8),
    @(132, 45874, "Mattheo Haon", 70, 8, 8, 0, "", 8),
    @(133, 45874, "Hedi Nasri", 70, 7, 5, 0, "", 6),
    @(134, 45874, "Wael Fareh", 70, 7, 7, 1, "Genou", 7),
    @(135, 45874, "Ilan Ihaddadene", 70, 8, 6, 2, "Ischio", 5),
    @(136, 45874, "Karahali Souaré", 70, 6, 6, 7, "Quadri", 8),
    @(137, 45874, "Amine Taiar", 70, 7, 8, 5, "Ischio", 8),
    @(138, 45874, "Naim Dhib", 70, 9, 10, 3, "Adducteur", 7),
    @(139, 45874, "Yoan Zouma", 70, 6, 7, 2, "Adducteur / ischio", 1),
    @(140, 45874, "Sofiane Belle", 70, 6, 4, 1, "Genou", 3),
    @(141, 45874, "Yanis Berrached", 70, 8, 8, 0, "", 8)
)

foreach ($row in $newRows) {
    $r      = $row[0]
    $date   = $row[1]
    $name   = $row[2]
    $volume = $row[3]
    $inten  = $row[4]
    $charge = $row[5]
    $fatig  = $row[6]
    $loc    = $row[7]
    $plais  = $row[8]

    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $volume
    $ws.Cells.Item($r, 4).Value = $inten
    $ws.Cells.Item($r, 5).Value = $charge
    $ws.Cells.Item($r, 6).Value = $fatig
    if ($loc -ne "") {
        $ws.Cells.Item($r, 7).Value = $loc
    } else {
        $ws.Range("G2").Copy()
        $ws.Cells.Item($r, 7).PasteSpecial(-4122)
        $excel.CutCopyMode = $false
    }
    $ws.Cells.Item($r, 8).Value = $plais
}

# 3) Formula column (I): mirrors how Excel splits the autofilled shared formula
$ws.Range("I126:I130").Formula = "=C126*D126"
$ws.Range("I131:I141").Formula = "=C131*D131"

# 4) Update the view state: scroll position + active selection
$ws.Range("K136").Select()

$wb.Application.Calculate()
